$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "26.072.90"
$ws.Range("D3").Value = "1.646.89"
$ws.Range("E3").Value = "  +0.55%  "
$ws.Range("E4").Value = "  +0.52%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "216.17"
$ws.Range("E5").Value = "  +0.65%  "
$ws.Range("E6").Value = "  +0.35%  "
$ws.Range("E7").Value = "  +0.61%  "
$ws.Range("E8").Value = "  +0.79%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.256"
$ws.Range("E9").Value = "  +0.62%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "19.62"
$ws.Range("E10").Value = "  +0.09%  "
$ws.Range("B12").Value = "WrappedEther"
$ws.Range("C12").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D12").Value = "1.691.46"
$ws.Range("E12").Value = "  +2.97%  "
$ws.Range("B13").Value = "Polkadot"
$ws.Range("C13").Value = "https://coinranking.com/coin/25W7FG7om+polkadot-dot"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "4.27"
$ws.Range("E13").Value = "  +0.74%  "
$ws.Range("E14").Value = "  +0.20%  "
$ws.Range("B15").Value = "Litecoin"
$ws.Range("C15").Value = "https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "63.58"
$ws.Range("E15").Value = "  +1.67%  "
$ws.Range("B16").Value = "ShibaInu"
$ws.Range("C16").Value = "https://coinranking.com/coin/xz24e0BjL+shibainu-shib"
$ws.Range("D16").Value = "0.0₃0764"
$ws.Range("E16").Value = "  +1.01%  "
$ws.Range("D17").Value = "26.080.69"
$ws.Range("E18").Value = "  +0.59%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "194.48"
$ws.Range("E19").Value = "  +0.53%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "4.37"
$ws.Range("E20").Value = "  -0.23%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "9.95"
$ws.Range("E21").Value = "  +0.34%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "6.22"
$ws.Range("E22").Value = "  -0.90%  "
$ws.Range("E23").Value = "  +4.85%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "1.80"
$ws.Range("E24").Value = "  -0.03%  "
$ws.Range("E25").Value = "  +0.60%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "143.86"
$ws.Range("E26").Value = "  -0.19%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "6.90"
$ws.Range("E27").Value = "  +0.69%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "15.53"
$ws.Range("E28").Value = "  +0.52%  "
$ws.Range("E29").Value = "  +0.64%  "
$ws.Range("E30").Value = "  -0.72%  "
$ws.Range("E31").Value = "  +1.56%  "
$ws.Range("E32").Value = "  -0.16%  "
$ws.Range("E33").Value = "  -0.18%  "
$ws.Range("E34").Value = "  +1.51%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.907"
$ws.Range("E35").Value = "  +0.49%  "
$ws.Range("D36").Value = "1.132.63"
$ws.Range("E36").Value = "  -0.39%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.540"
$ws.Range("E37").Value = "  -0.90%  "
$ws.Range("E38").Value = "  +0.42%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.0158"
$ws.Range("E39").Value = "  +0.60%  "
$ws.Range("E40").Value = "  +0.84%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "99.06"
$ws.Range("E41").Value = "  -0.32%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.799"
$ws.Range("E42").Value = "  +0.27%  "
$ws.Range("E43").Value = "  +2.15%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "56.66"
$ws.Range("E44").Value = "  +0.29%  "
$ws.Range("E45").Value = "  +2.90%  "
$ws.Range("E46").Value = "  -1.26%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "7.78"
$ws.Range("E47").Value = "  +1.61%  "
$ws.Range("E48").Value = "  +0.01%  "
$ws.Range("E49").Value = "  +0.33%  "
$ws.Range("E50").Value = "  -1.02%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "1.19"
$ws.Range("E51").Value = "  +3.01%  "
